$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new blank rows before the current row 14 so the old rows
#    14-18 (blank spacer rows, total row, and discount notes) shift down to
#    16-20. The two inserted rows inherit formatting (styles) from the row
#    above, matching rows 14/15 in the target file.
# ---------------------------------------------------------------------------
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).Insert()

# ---------------------------------------------------------------------------
# 2. Populate the four new component rows (10-13) with the newly found SMD
#    resistors and radial capacitors.
# ---------------------------------------------------------------------------

# Row 10: Vishay Dale CRCW2010100RFKEFHP (100 ohm SMD resistor)
$ws.Range("A10").Value = "Vishay Dale CRCW2010100RFKEFHP"
$ws.Range("B10").Value = "Solderable size 2010 100 ohm SMD"
$ws.Range("C10").Value = 0.58
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 18
$ws.Range("F10").Formula = "=C10*D10"
$ws.Range("G10").Value = "Digi-key"
$ws.Range("H10").Value = "https://www.digikey.com/product-detail/en/vishay-dale/CRCW2010100RFKEFHP/541-100PCT-ND/2222676"

# Row 11: Vishay Dale CRCW201010R0FKEF (10 ohm SMD resistor)
$ws.Range("A11").Value = "Vishay Dale CRCW201010R0FKEF"
$ws.Range("B11").Value = "Solderable size 2010 10 ohm SMD"
$ws.Range("C11").Value = 0.29
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 6
$ws.Range("F11").Formula = "=C11*D11"
$ws.Range("G11").Value = "Digi-key"
$ws.Range("H11").Value = "https://www.digikey.com/product-detail/en/vishay-dale/CRCW201010R0FKEF/541-10.0ACCT-ND/1179050"

# Row 12: WIMA FKP0D001000B00JSSD (100 pF film cap)
$ws.Range("A12").Value = "WIMA FKP0D001000B00JSSD"
$ws.Range("B12").Value = "Radial 100 pF Film Cap "
$ws.Range("C12").Value = 0.7
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = 9
$ws.Range("F12").Formula = "=C12*D12"
$ws.Range("G12").Value = "Digi-key"
$ws.Range("H12").Value = "https://www.digikey.com/product-detail/en/wima/FKP0D001000B00JSSD/1928-1039-ND/9370034"

# Row 13: Nichicon UVK2GR47MED1TD (0.1 uF electrolytic cap)
$ws.Range("A13").Value = "Nichicon UVK2GR47MED1TD "
$ws.Range("B13").Value = "Radial 0.1 uF Electrolytic Cap"
$ws.Range("C13").Value = 0.06
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 3
$ws.Range("F13").Formula = "=C13*D13"
$ws.Range("G13").Value = "Digi-key"
$ws.Range("H13").Value = "https://www.digikey.com/product-detail/en/nichicon/UVK2GR47MED1TD/493-12648-3-ND/4328849"

# ---------------------------------------------------------------------------
# 3. Apply the blue "link" look (same style used on the other Link column
#    cells that don't carry a real hyperlink object, e.g. H4/H6/H7/H8 ->
#    here H2/H3/H5/H9's formatting) to the four new Link cells.
# ---------------------------------------------------------------------------
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H10:H13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Extend the running total formula to include the four new rows.
# ---------------------------------------------------------------------------
$ws.Range("F17").Formula = "=F2+F3+F4+F5+F6+F7+F8+F9+F10+F11+F12+F13"

# ---------------------------------------------------------------------------
# 5. Widen column A to fit the new, longer component names.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 32.75

# ---------------------------------------------------------------------------
# 6. Update the active cell/selection to reflect where editing left off.
# ---------------------------------------------------------------------------
$ws.Range("B19").Select() | Out-Null
